# The workbook's "data" sheet has a block of cells Y2:Y20 that currently
# hold the shared string "saudacoes"; the author re-tagged that whole
# column of rows to "feedback_positivo". Re-apply the same edit here.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Y2:Y20").Value = "feedback_positivo"

# The author's session also left the sheet scrolled so column R is the
# first visible column (topLeftCell went from O1 to R1), with the same
# Y2:Y20 selection as before. Scroll the active window to match.
$ws.Range("Y2:Y20").Select()
$excel.ActiveWindow.ScrollColumn = 18
$excel.ActiveWindow.ScrollRow = 1
